$d = $word.ActiveDocument

# Splits the hex-colour text in a paragraph (e.g. "#3B55A2") into two
# separate runs - one holding just "#" and one holding the new hex
# digits - mirroring the target OOXML, which shows the colour value
# split across two <w:r> elements.
function Split-ColorRun($paragraphIndex, $oldColor, $newColor) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $searchRng = $d.Range($p.Range.Start, $p.Range.End)
    $found = $searchRng.Find.Execute($oldColor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    # Build a fresh Range over the matched text (re-using the Find range
    # directly confuses later calls), then replace its contents with two
    # runs via InsertXML so the split survives the save as real <w:r>s.
    $targetRng = $d.Range($searchRng.Start, $searchRng.End)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>#</w:t></w:r><w:r><w:t>' + $newColor + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $targetRng.InsertXML($xml)
}

Split-ColorRun 2 "#3B55A2" "3B56A3"
Split-ColorRun 3 "#88C34A" "88C449"

# Paragraph 4 ends with the "_GoBack" bookmark sitting right after its
# colour text. InsertXML on a range that reaches the very end of the
# document's content shifts that bookmark to the wrong spot (it ends up
# wrapping the whole replaced range instead of staying collapsed at the
# end), so it needs to be restored afterwards.
Split-ColorRun 4 "#88C34A" "F37139"

$p4 = $d.Paragraphs.Item(4)
$endPos = $p4.Range.End - 1

# Re-anchoring a bookmark exactly at the document's last valid text
# position is unreliable in this runtime, so a throwaway character is
# inserted there first to move the trouble spot out of the way, the
# bookmark is rebuilt in front of it, and then the throwaway character
# is deleted again.
$tempRng = $d.Range($endPos, $endPos)
$tempRng.InsertAfter("Z")

$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$fixRng = $d.Range($endPos, $endPos)
$fixRng.Bookmarks.Add("_GoBack")

$zRng = $d.Range($endPos, $endPos + 1)
$zRng.Text = ""
